$wb = $excel.ActiveWorkbook

# Create the new "UK" sheet by copying the existing "Poland" sheet (same
# layout/styles/merged cells) and placing it after the last sheet.
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# Insert a new row above the old row 9 ("PR1D2") for the "GMPIM" entry,
# copying the formatting from row 8 so the new cell keeps the same style.
$uk.Rows.Item(9).Insert()
$uk.Range("A8").Copy($uk.Range("A9"))
$uk.Range("A9").Value = "GMPIM"

# Update the user-story/NGC reference and market name for the UK sheet.
$uk.Range("B4").Value = "NGC-2741/T3365/T3366/T3364"
$uk.Range("B2").Value = "UK Market"

# Match the selection/active cell shown in the saved sheet.
$uk.Range("A9").Select()
